$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 content (set in this order so shared-string table indices line up
# with the target: "Kerckhoff", "Kerckhoffs' Prinzip", then the long C9 text)
$ws.Range("A10").Value = "Kerckhoff"
$ws.Range("B10").Value = "Kerckhoffs’ Prinzip"
$ws.Range("C9").Value = "Das Kerckhoffs’sche Prinzip oder Kerckhoffs’ Maxime ist ein im Jahr 1883 von Auguste Kerckhoffs formulierter Grundsatz der modernen Kryptographie, welcher besagt, dass die Sicherheit eines (symmetrischen) Verschlüsselungsverfahrens auf der Geheimhaltung des Schlüssels beruht anstatt auf der Geheimhaltung des Verschlüsselungsalgorithmus. "

# Match formatting of equivalent existing rows
$ws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

# Row heights grew to fit the newly-wrapped long text / new row
$ws.Rows.Item(9).RowHeight = 105.75
$ws.Rows.Item(10).RowHeight = 26.25
$ws.Rows.Item(2).RowHeight = 25.5

# Cursor/selection ended up on B5
[void]$ws.Range("B5").Select()
